$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml): update F3 and F4
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 1187
$wsExhibit.Range("F4").Value = 2655

# Sheet "全部类型" (sheet4.xml): update F5 and F6
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 1187
$wsAll.Range("F6").Value = 2655
